# Update database: drop oldest "6 ماهه منتهی به 1399/06" period, shift every
# period/date/value one column to the left, and append the new
# "12 ماهه منتهی به 1401/12" period (published 1402-02-27) as the newest column (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: period headers (columns D..M) ---
$periods = @(
    "9 ماهه منتهی به 1399/09",
    "12 ماهه منتهی به 1399/12",
    "3 ماهه منتهی به 1400/03",
    "6 ماهه منتهی به 1400/06",
    "9 ماهه منتهی به 1400/09",
    "12 ماهه منتهی به 1400/12",
    "3 ماهه منتهی به 1401/03",
    "6 ماهه منتهی به 1401/06",
    "9 ماهه منتهی به 1401/09",
    "12 ماهه منتهی به 1401/12"
)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $periods[$i]
}

# --- Row 9: publish dates (columns D..M) ---
$dates = @(
    "1400-10-29 (2)",
    "1401-03-24 (9)",
    "1401-04-30 (3)",
    "1401-08-29 (4)",
    "1401-10-28 (2)",
    "1402-02-27 (7)",
    "1401-04-30",
    "1401-08-29 (2)",
    "1401-10-28",
    "1402-02-27"
)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $dates[$i]
}

# --- Data rows 11-27 (columns D..M), each shifted one period to the left,
#     with the new 12-month-cumulative figure landing in column M ---
$data = @{
    11 = @(734347, 1008240, 217182, 549295, 1019462, 1632271, 590461, 1370075, 2131043, 2999738)
    12 = @(-557701, -747453, -139110, -372590, -719151, -1139548, -460541, -1110304, -1746984, -2387738)
    13 = @(176646, 260787, 78072, 176705, 300311, 492723, 129920, 259771, 384059, 612000)
    14 = @(-24115, -33021, -14060, -30516, -42490, -60712, -14340, -37350, -60526, -95353)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(4150, 4265, 999, 4172, 2272, 3929, 2795, 6124, 8788, 34059)
    17 = @(156681, 232031, 65011, 150361, 260093, 435940, 118375, 228545, 332321, 550706)
    18 = @(-28028, -44626, -28433, -50069, -85677, -168851, -30208, -76790, -130613, -262925)
    19 = @(3803, 1477, 671, -1196, -34485, 14983, 6848, 14266, 23603, 29352)
    20 = @(132456, 188882, 37249, 99096, 139931, 282072, 95015, 166021, 225311, 317133)
    21 = @(-22514, -32796, -6401, -19316, -23455, -57039, -13575, -32436, -33450, -40052)
    22 = @(109942, 156086, 30848, 79780, 116476, 225033, 81440, 133585, 191861, 277081)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(109942, 156086, 30848, 79780, 116476, 225033, 81440, 133585, 191861, 277081)
    25 = @(142, 201, 40, 103, 150, 290, 105, 172, 248, 358)
    26 = @(775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000)
    27 = @(142, 201, 40, 103, 150, 290, 105, 172, 248, 358)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
